# Populate the USB sheet's histogram data, switch the active tab to USB,
# and update the selection on both sheets accordingly.

$wb = $excel.ActiveWorkbook

$ssd = $wb.Worksheets.Item("SSD")
$usb = $wb.Worksheets.Item("USB")

# Histogram data (Time in MS / Rate (bytes/s)) for the USB sheet.
$data = @(
    @(1416, 74),
    @(4993, 21),
    @(7489, 14),
    @(5518, 19),
    @(5242, 20),
    @(7489, 14),
    @(6553, 16),
    @(9532, 11),
    @(5242, 20),
    @(4369, 24)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $usb.Cells.Item($row, 2).Value = $data[$i][0]
    $usb.Cells.Item($row, 3).Value = $data[$i][1]
}

# Update the selections on each sheet.
$ssd.Range("E8").Select() | Out-Null
$usb.Range("B9").Select() | Out-Null

# Make USB the active sheet/tab.
$usb.Activate()
